$d = $word.ActiveDocument

# --- 1. Remove the "TILBAGE:" paragraph and the blank paragraph right after it ---
$findRange = $d.Content
$found = $findRange.Find.Execute("TILBAGE:", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if ($found -and $findRange.Find.Found) {
    $tilbagePara = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -le $findRange.Start -and $p.Range.End -ge $findRange.End) {
            $tilbagePara = $p
            break
        }
    }
    if ($tilbagePara -ne $null) {
        $nextPara = $tilbagePara.Next()
        $deleteRange = $d.Range($tilbagePara.Range.Start, $nextPara.Range.End)
        $deleteRange.Delete()
    }
}

# --- 2. Add a new paragraph right after "Hvor meget skal server-side renders?" ---
$findRange2 = $d.Content
$found2 = $findRange2.Find.Execute("Hvor meget skal server-side renders?", $true, $false, $false, `
                                    $false, $false, $true, 1, $false, "", 0)
if ($found2 -and $findRange2.Find.Found) {
    $serverPara = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -le $findRange2.Start -and $p.Range.End -ge $findRange2.End) {
            $serverPara = $p
            break
        }
    }
    if ($serverPara -ne $null) {
        $serverPara.Range.InsertParagraphAfter()
        $newPara = $serverPara.Next()
        $oe = [char]248
        $newPara.Range.Text = "Jeg tror jeg vil lave det meste back-end f" + $oe + "rst og bekymre mig om front-end bagefter"
    }
}
